$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2 already uses the shared style index 1 in the original workbook; mutate its
# number format to Text ("@" -> numFmtId 49), then apply the same format to the
# other cells in the row so they reuse that same style record.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("A2").NumberFormat = "@"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("D2").NumberFormat = "@"

$ws.Range("A2").Value = "21 V"
$ws.Range("B2").Value = "3 A"
$ws.Range("D2").Value = "322"
$ws.Range("C2").Value = "100 %"

$ws.Range("C3").Select()
